# 2.1.1.1e — add 2023 data column (Q) and adjust row heights / selection.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. New column Q: copy formatting from the corresponding P cell, then
#    set the 2023 values (row 3 stays blank, just inherits the border
#    style used across the rest of that separator row).
# ---------------------------------------------------------------------
$newValues = @{
    4  = 2023
    5  = 74.605426356589135
    6  = 118.8
    7  = 71.61643835616438
    8  = 95.703125
    9  = 113.91018619934282
    10 = 108.21501014198785
    11 = 165.26684164479443
    12 = 48.504446240905416
    13 = 97.361348644026393
    14 = 52.747252747252752
}

# Row 3 only needs the formatting of the thick-bottom border copied over.
$ws.Range("P3").Copy()
$ws.Range("Q3").PasteSpecial(-4122)

foreach ($row in 4..14) {
    $srcCell = $ws.Cells.Item($row, 16)   # column P
    $dstCell = $ws.Cells.Item($row, 17)   # column Q
    $srcCell.Copy()
    $dstCell.PasteSpecial(-4122)
    $dstCell.Value = $newValues[$row]
}

# ---------------------------------------------------------------------
# 2. Row heights (rows 4-14 gain/adjust an explicit custom height).
# ---------------------------------------------------------------------
$ws.Rows.Item(4).RowHeight = 16.5
$ws.Rows.Item(5).RowHeight = 27
$ws.Rows.Item(6).RowHeight = 24.75
foreach ($row in 7..14) {
    $ws.Rows.Item($row).RowHeight = 16.5
}

# ---------------------------------------------------------------------
# 3. Reset the saved selection back to the default top-left cell.
# ---------------------------------------------------------------------
$ws.Range("A1").Select() | Out-Null
